$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.340.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.57%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.825.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.75%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -3.06%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'314.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.48%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -2.85%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4294"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -2.54%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3696"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -3.27%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07243"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.72%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.8642"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.66%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'21.13"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.21%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.826.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.90%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'6.664"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.33%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.344"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.46%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.07084"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.72%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'87.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.70%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.004"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -3.29%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008870"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.75%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -2.91%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'15.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.35%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'27.358.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.59%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.157"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.81%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'10.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.90%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.051.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.90%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.009"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -3.33%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'153.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -3.75%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'18.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.00%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.139"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +7.15%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'5.285"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.45%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'116.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.14%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.08837"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.13%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.201"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.18%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.7649"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.25%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.490"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.50%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.850"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -6.31%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -3.07%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.122"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.99%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.01958"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.70%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05261"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.21%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.882"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.74%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'7.109"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.07%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1679"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.27%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.5056"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.17%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'8.648"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.91%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'10.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.36%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'106.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.36%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.4716"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.56%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.06425"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.17%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.003"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.16%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.666"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.01%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.816"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -4.08%  "
$ws.Range("E51").Style = "Normal"
